# Fixed build errors caused by Resharper missing some replacements when
# doing a rename refactor, and weird names chosen by Resharper.
#
# Concretely: a couple of long "zzzz..." test/placeholder strings were
# extended by a few characters, and several stale numeric sample values
# (dates/amounts) were bumped to their corrected figures.

$wb = $excel.ActiveWorkbook

$budgetOut  = $wb.Worksheets.Item("Budget Out")
$testRecord = $wb.Worksheets.Item("TestRecord")
$expectedOut = $wb.Worksheets.Item("Expected Out")

# --- Shared/placeholder text tweaks ---
$budgetOut.Range("F9").Value = "Description007zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"
$testRecord.Range("E10").Value = "some test textzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# --- Corrected numeric values ---
$expectedOut.Range("B9").Value = 1351.76
$expectedOut.Range("B11").Value = 431.62

$testRecord.Range("A10").Value = 43266
$testRecord.Range("B10").Value = 124.74

$budgetOut.Range("C9").Value = 93.82

# 'Expected Out'!B1 is =SUM(B2:B295); it recalculates automatically once
# B9/B11 change, so no direct write is needed there.
